$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("2000s")

$ws.Range("A16").Value = "La vereda de la puerta de atrás"
$ws.Range("B16").Value = "Extremoduro"
$ws.Range("C16").Value = "Yo, Minoría Absoluta"
$ws.Range("D16").Value = 2002
$ws.Range("E16").Value = "spotify:track:4kJP8Z888wREJ8bRMWNMuk"
